# Update the "取得日時" (acquired datetime) column on the "ランサーズ" sheet.
# All existing data rows (2-10) get their timestamp bumped from
# 2025-12-18 12:38:54 to 2025-12-18 12:51:43 (new scrape pass at that time).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-18 12:38:54"
$newTimestamp = "2025-12-18 12:51:43"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
